$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "High Scores" -> "High Scores (1 - 5)"
#    Append a new, separate run containing " (1 - 5)" (en dash) right
#    after the existing "High Scores" run, inside the same paragraph.
# ---------------------------------------------------------------------
$hsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "High Scores") {
        $hsPara = $p
        break
    }
}

$hsRange = $hsPara.Range
$hsRange.End = $hsRange.End - 1
$insertStart = $hsRange.End
$suffixText = " (1 " + [char]0x2013 + " 5)"
$hsRange.InsertAfter($suffixText)

$suffixRange = $d.Range($insertStart, $insertStart + $suffixText.Length)
# Round-tripping FormattedText forces the new text to live in its own
# <w:r> run instead of being merged back into the preceding run.
$suffixRange.FormattedText = $suffixRange.FormattedText

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the very end of the document
#    (after "Game Name") to right after the colon in the
#    "users-charities:" paragraph.
# ---------------------------------------------------------------------
$ucPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "users-charities:") {
        $ucPara = $p
        break
    }
}

$ucRange = $ucPara.Range
$ucRange.End = $ucRange.End - 1
$markerStart = $ucRange.End
$marker = "@@BOOKMARK_MARKER@@"
$ucRange.InsertAfter($marker)

$markerRange = $d.Content
$found = $markerRange.Find.Execute($marker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $markerRange)

$markerRange2 = $d.Content
$found2 = $markerRange2.Find.Execute($marker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerRange2.Delete()
